$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows: first_name, last_name, height, weight, date_of_birth (as a date serial)
$data = @(
    @("Adam",   "Smith",  11.5, 180.5, 29545),
    @("Janice", "Harper", 10.4, 160.3, 33727),
    @("Joshua", "Kooler", 12,   190.6, 33669),
    @("Amelia", "Sholer", 10.5, 160.4, 32998),
    @("Peter",  "Jiang",  11.5, 140.5, 34853)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value  = $row[0]
    $ws.Cells.Item($r, 2).Value  = $row[1]
    $ws.Cells.Item($r, 3).Value  = $row[2]
    $ws.Cells.Item($r, 4).Value  = $row[3]
    $ws.Cells.Item($r, 5).Value2 = $row[4]
    $r++
}

# Format the date_of_birth column as dates (built-in mm-dd-yy format), applying the
# same style to every cell (set once, then copy the format to the rest).
$firstDateCell = $ws.Cells.Item(2, 5)
$firstDateCell.NumberFormat = "mm-dd-yy"
$firstDateCell.Copy()
$ws.Range("E3:E6").PasteSpecial(-4122)  # xlPasteFormats
